$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "93.799.33"
Set-TextValue "E2" "  -1.34%  "
Set-TextValue "D3" "3.323.92"
Set-TextValue "E3" "  -3.23%  "
Set-TextValue "E4" "  -0.06%  "
Set-TextValue "D5" "230.64"
Set-TextValue "E5" "  -3.26%  "
Set-TextValue "D6" "617.34"
Set-TextValue "E6" "  -3.55%  "
Set-TextValue "D7" "1.38"
Set-TextValue "E7" "  -3.64%  "
Set-TextValue "D8" "0.387"
Set-TextValue "E8" "  -3.15%  "
Set-TextValue "E9" "  -0.06%  "
Set-TextValue "D10" "0.926"
Set-TextValue "E10" "  -5.98%  "
Set-TextValue "D11" "3.323.67"
Set-TextValue "E11" "  -3.16%  "
Set-TextValue "D12" "42.05"
Set-TextValue "E12" "  +1.72%  "
Set-TextValue "E13" "  -1.85%  "
Set-TextValue "B14" "WrappedBTC"
Set-TextValue "C14" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D14" "93.588.18"
Set-TextValue "E14" "  -1.31%  "
Set-TextValue "B15" "Toncoin"
Set-TextValue "C15" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D15" "5.94"
Set-TextValue "E15" "  -2.24%  "
Set-TextValue "D16" "3.950.55"
Set-TextValue "E16" "  -3.19%  "
Set-TextValue "E17" "  -4.68%  "
Set-TextValue "D18" "8.09"
Set-TextValue "E18" "  -3.42%  "
Set-TextValue "D19" "3.325.50"
Set-TextValue "E19" "  -3.44%  "
Set-TextValue "D20" "17.16"
Set-TextValue "E20" "  -3.92%  "
Set-TextValue "E21" "  -5.20%  "
Set-TextValue "E22" "  +9.91%  "
Set-TextValue "D23" "494.83"
Set-TextValue "E23" "  -0.97%  "
Set-TextValue "D24" "0.447"
Set-TextValue "E24" "  -12.25%  "
Set-TextValue "E25" "  -4.53%  "
Set-TextValue "D26" "6.14"
Set-TextValue "E26" "  -6.13%  "
Set-TextValue "D27" "91.50"
Set-TextValue "E27" "  +0.20%  "
Set-TextValue "D28" "11.70"
Set-TextValue "E28" "  -2.24%  "
Set-TextValue "D29" "3.505.97"
Set-TextValue "E29" "  -3.20%  "
Set-TextValue "E30" "  -0.12%  "
Set-TextValue "D31" "11.02"
Set-TextValue "E31" "  -4.99%  "
Set-TextValue "E32" "  +1.35%  "
Set-TextValue "D33" "2.60"
Set-TextValue "E33" "  -4.15%  "
Set-TextValue "D34" "0.996"
Set-TextValue "E34" "  -0.06%  "
Set-TextValue "D35" "0.174"
Set-TextValue "E35" "  -4.58%  "
Set-TextValue "D36" "28.26"
Set-TextValue "E36" "  -7.77%  "
Set-TextValue "D37" "0.528"
Set-TextValue "E37" "  -6.15%  "
Set-TextValue "D38" "529.05"
Set-TextValue "E38" "  +3.63%  "
Set-TextValue "D39" "7.38"
Set-TextValue "E39" "  -3.89%  "
Set-TextValue "E40" "  +0.09%  "
Set-TextValue "E41" "  -1.23%  "
Set-TextValue "E42" "  -5.22%  "
Set-TextValue "D43" "0.859"
Set-TextValue "E43" "  -5.00%  "
Set-TextValue "B44" "WhiteBITCoin"
Set-TextValue "C44" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D44" "24.04"
Set-TextValue "E44" "  -0.37%  "
Set-TextValue "B45" "MantraDAO"
Set-TextValue "C45" "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
Set-TextValue "D45" "3.66"
Set-TextValue "E45" "  +4.32%  "
Set-TextValue "D46" "0.0414"
Set-TextValue "E46" "  +0.46%  "
Set-TextValue "D47" "1.67"
Set-TextValue "E47" "  -0.60%  "
Set-TextValue "D48" "5.38"
Set-TextValue "E48" "  -2.20%  "
Set-TextValue "D49" "52.87"
Set-TextValue "E49" "  -1.05%  "
Set-TextValue "D50" "2.11"
Set-TextValue "E50" "  -1.03%  "
Set-TextValue "D51" "7.93"
Set-TextValue "E51" "  +0.13%  "
